$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1052975711"
$ws.Range("D16").Value = "CATHERINE MARCELA CRUZ CEBALLOS"
$ws.Range("E16").Value = "2202"
$ws.Range("F16").Value = 25439
$ws.Range("G16").Value = 908526

$ws.Range("C17").Value = "1052975711"
$ws.Range("D17").Value = "CATHERINE MARCELA CRUZ CEBALLOS"
$ws.Range("E17").Value = "2201"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526

$ws.Range("C18").Value = "1052975711"
$ws.Range("D18").Value = "CATHERINE MARCELA CRUZ CEBALLOS"
$ws.Range("E18").Value = "2112"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526

$ws.Range("C19").Value = "1052975711"
$ws.Range("D19").Value = "CATHERINE MARCELA CRUZ CEBALLOS"
$ws.Range("E19").Value = "2111"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

$ws.Range("C20").Value = "1052975711"
$ws.Range("D20").Value = "CATHERINE MARCELA CRUZ CEBALLOS"
$ws.Range("E20").Value = "2110"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 908526

$ws.Range("C21").Value = "73238615"
$ws.Range("D21").Value = "LESTER MANUEL SABAYE VANEGAS"
$ws.Range("E21").Value = "2202"
$ws.Range("F21").Value = 25439
$ws.Range("G21").Value = 908526

$ws.Range("C22").Value = "73238615"
$ws.Range("D22").Value = "LESTER MANUEL SABAYE VANEGAS"
$ws.Range("E22").Value = "2201"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

$ws.Range("C23").Value = "73238615"
$ws.Range("D23").Value = "LESTER MANUEL SABAYE VANEGAS"
$ws.Range("E23").Value = "2112"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 908526

$ws.Range("C24").Value = "73238615"
$ws.Range("D24").Value = "LESTER MANUEL SABAYE VANEGAS"
$ws.Range("E24").Value = "2111"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 908526

$ws.Range("C25").Value = "73238615"
$ws.Range("D25").Value = "LESTER MANUEL SABAYE VANEGAS"
$ws.Range("E25").Value = "2110"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 908526

$ws.Range("C26").Value = "1052996035"
$ws.Range("D26").Value = "JULIETH CANDELARIA CASTILLO PEREZ"
$ws.Range("E26").Value = "2202"
$ws.Range("F26").Value = 25439
$ws.Range("G26").Value = 908526

$ws.Range("C27").Value = "1052996035"
$ws.Range("D27").Value = "JULIETH CANDELARIA CASTILLO PEREZ"
$ws.Range("E27").Value = "2201"
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = 908526

$ws.Range("C28").Value = "1052996035"
$ws.Range("D28").Value = "JULIETH CANDELARIA CASTILLO PEREZ"
$ws.Range("E28").Value = "2112"
$ws.Range("F28").Value = 36341
$ws.Range("G28").Value = 908526

$ws.Range("C29").Value = "1052996035"
$ws.Range("D29").Value = "JULIETH CANDELARIA CASTILLO PEREZ"
$ws.Range("E29").Value = "2111"
$ws.Range("F29").Value = 36341
$ws.Range("G29").Value = 908526

$ws.Range("C30").Value = "1052996035"
$ws.Range("D30").Value = "JULIETH CANDELARIA CASTILLO PEREZ"
$ws.Range("E30").Value = "2110"
$ws.Range("F30").Value = 36341
$ws.Range("G30").Value = 908526

$ws.Range("C31").Value = "1005646781"
$ws.Range("D31").Value = "WILMER HERNANDO QUIROZ CATAÑO"
$ws.Range("E31").Value = "2202"
$ws.Range("F31").Value = 25439
$ws.Range("G31").Value = 908528

$ws.Range("C32").Value = "1005646781"
$ws.Range("D32").Value = "WILMER HERNANDO QUIROZ CATAÑO"
$ws.Range("E32").Value = "2201"
$ws.Range("F32").Value = 36341
$ws.Range("G32").Value = 908528

$ws.Range("C33").Value = "1005646781"
$ws.Range("D33").Value = "WILMER HERNANDO QUIROZ CATAÑO"
$ws.Range("E33").Value = "2112"
$ws.Range("F33").Value = 36341
$ws.Range("G33").Value = 908528

$ws.Range("C34").Value = "1005646781"
$ws.Range("D34").Value = "WILMER HERNANDO QUIROZ CATAÑO"
$ws.Range("E34").Value = "2111"
$ws.Range("F34").Value = 36341
$ws.Range("G34").Value = 908528

$ws.Range("C35").Value = "1005646781"
$ws.Range("D35").Value = "WILMER HERNANDO QUIROZ CATAÑO"
$ws.Range("E35").Value = "2110"
$ws.Range("F35").Value = 36341
$ws.Range("G35").Value = 908528

$ws.Range("C36").Value = "1005646185"
$ws.Range("D36").Value = "JAIR JOSE OSORIO GARCIA"
$ws.Range("E36").Value = "2202"
$ws.Range("F36").Value = 25439
$ws.Range("G36").Value = 908526

$ws.Range("C37").Value = "1005646185"
$ws.Range("D37").Value = "JAIR JOSE OSORIO GARCIA"
$ws.Range("E37").Value = "2201"
$ws.Range("F37").Value = 36341
$ws.Range("G37").Value = 908526

$ws.Range("C38").Value = "1005646185"
$ws.Range("D38").Value = "JAIR JOSE OSORIO GARCIA"
$ws.Range("E38").Value = "2112"
$ws.Range("F38").Value = 36341
$ws.Range("G38").Value = 908526

$ws.Range("C39").Value = "1005646185"
$ws.Range("D39").Value = "JAIR JOSE OSORIO GARCIA"
$ws.Range("E39").Value = "2111"
$ws.Range("F39").Value = 36341
$ws.Range("G39").Value = 908526

$ws.Range("C40").Value = "1005646185"
$ws.Range("D40").Value = "JAIR JOSE OSORIO GARCIA"
$ws.Range("E40").Value = "2110"
$ws.Range("F40").Value = 36341
$ws.Range("G40").Value = 908526

$ws.Range("C41").Value = "1052986417"
$ws.Range("D41").Value = "JORGE ISAAC NOYA PRASCA"
$ws.Range("E41").Value = "2202"
$ws.Range("F41").Value = 25439
$ws.Range("G41").Value = 908528

$ws.Range("C42").Value = "1052986417"
$ws.Range("D42").Value = "JORGE ISAAC NOYA PRASCA"
$ws.Range("E42").Value = "2201"
$ws.Range("F42").Value = 36341
$ws.Range("G42").Value = 908528

$ws.Range("C43").Value = "1052986417"
$ws.Range("D43").Value = "JORGE ISAAC NOYA PRASCA"
$ws.Range("E43").Value = "2112"
$ws.Range("F43").Value = 36341
$ws.Range("G43").Value = 908528

$ws.Range("C44").Value = "1052986417"
$ws.Range("D44").Value = "JORGE ISAAC NOYA PRASCA"
$ws.Range("E44").Value = "2111"
$ws.Range("F44").Value = 36341
$ws.Range("G44").Value = 908528

$ws.Range("C45").Value = "1052986417"
$ws.Range("D45").Value = "JORGE ISAAC NOYA PRASCA"
$ws.Range("E45").Value = "2110"
$ws.Range("F45").Value = 36341
$ws.Range("G45").Value = 908528

$ws.Range("C46").Value = "73238056"
$ws.Range("D46").Value = "JORGE ELIECER LUNA RODELO"
$ws.Range("E46").Value = "2202"
$ws.Range("F46").Value = 25439
$ws.Range("G46").Value = 908526

$ws.Range("C47").Value = "73238056"
$ws.Range("D47").Value = "JORGE ELIECER LUNA RODELO"
$ws.Range("E47").Value = "2201"
$ws.Range("F47").Value = 36341
$ws.Range("G47").Value = 908526

$ws.Range("C48").Value = "73238056"
$ws.Range("D48").Value = "JORGE ELIECER LUNA RODELO"
$ws.Range("E48").Value = "2112"
$ws.Range("F48").Value = 36341
$ws.Range("G48").Value = 908526

$ws.Range("C49").Value = "73238056"
$ws.Range("D49").Value = "JORGE ELIECER LUNA RODELO"
$ws.Range("E49").Value = "2111"
$ws.Range("F49").Value = 36341
$ws.Range("G49").Value = 908526

$ws.Range("C50").Value = "73238056"
$ws.Range("D50").Value = "JORGE ELIECER LUNA RODELO"
$ws.Range("E50").Value = "2110"
$ws.Range("F50").Value = 36341
$ws.Range("G50").Value = 908526

$ws.Range("C51").Value = "73243635"
$ws.Range("D51").Value = "RONALD ENRIQUE AVILA VANEGAS"
$ws.Range("E51").Value = "2202"
$ws.Range("F51").Value = 25439
$ws.Range("G51").Value = 908526

$ws.Range("C52").Value = "73243635"
$ws.Range("D52").Value = "RONALD ENRIQUE AVILA VANEGAS"
$ws.Range("E52").Value = "2201"
$ws.Range("F52").Value = 36341
$ws.Range("G52").Value = 908526

$ws.Range("C53").Value = "73243635"
$ws.Range("D53").Value = "RONALD ENRIQUE AVILA VANEGAS"
$ws.Range("E53").Value = "2112"
$ws.Range("F53").Value = 36341
$ws.Range("G53").Value = 908526

$ws.Range("C54").Value = "73243635"
$ws.Range("D54").Value = "RONALD ENRIQUE AVILA VANEGAS"
$ws.Range("E54").Value = "2111"
$ws.Range("F54").Value = 36341
$ws.Range("G54").Value = 908526

$ws.Range("C55").Value = "73243635"
$ws.Range("D55").Value = "RONALD ENRIQUE AVILA VANEGAS"
$ws.Range("E55").Value = "2110"
$ws.Range("F55").Value = 36341
$ws.Range("G55").Value = 908526

$ws.Range("C56").Value = "1052995157"
$ws.Range("D56").Value = "GILMAR ENRIQUE PARRA BASANTA"
$ws.Range("E56").Value = "2202"
$ws.Range("F56").Value = 25439
$ws.Range("G56").Value = 908526

$ws.Range("C57").Value = "1052995157"
$ws.Range("D57").Value = "GILMAR ENRIQUE PARRA BASANTA"
$ws.Range("E57").Value = "2201"
$ws.Range("F57").Value = 36341
$ws.Range("G57").Value = 908526

$ws.Range("C58").Value = "1052995157"
$ws.Range("D58").Value = "GILMAR ENRIQUE PARRA BASANTA"
$ws.Range("E58").Value = "2112"
$ws.Range("F58").Value = 36341
$ws.Range("G58").Value = 908526

$ws.Range("C59").Value = "1052995157"
$ws.Range("D59").Value = "GILMAR ENRIQUE PARRA BASANTA"
$ws.Range("E59").Value = "2111"
$ws.Range("F59").Value = 36341
$ws.Range("G59").Value = 908526

$ws.Range("C60").Value = "1052995157"
$ws.Range("D60").Value = "GILMAR ENRIQUE PARRA BASANTA"
$ws.Range("E60").Value = "2110"
$ws.Range("F60").Value = 36341
$ws.Range("G60").Value = 908526

$ws.Range("C61").Value = "33205661"
$ws.Range("D61").Value = "IRIS INES VANEGAS LUNA"
$ws.Range("E61").Value = "2202"
$ws.Range("F61").Value = 25439
$ws.Range("G61").Value = 908526

$ws.Range("C62").Value = "33205661"
$ws.Range("D62").Value = "IRIS INES VANEGAS LUNA"
$ws.Range("E62").Value = "2201"
$ws.Range("F62").Value = 36341
$ws.Range("G62").Value = 908526

$ws.Range("C63").Value = "33205661"
$ws.Range("D63").Value = "IRIS INES VANEGAS LUNA"
$ws.Range("E63").Value = "2112"
$ws.Range("F63").Value = 36341
$ws.Range("G63").Value = 908526

$ws.Range("C64").Value = "33205661"
$ws.Range("D64").Value = "IRIS INES VANEGAS LUNA"
$ws.Range("E64").Value = "2111"
$ws.Range("F64").Value = 36341
$ws.Range("G64").Value = 908526

$ws.Range("C65").Value = "33205661"
$ws.Range("D65").Value = "IRIS INES VANEGAS LUNA"
$ws.Range("E65").Value = "2110"
$ws.Range("F65").Value = 36341
$ws.Range("G65").Value = 908526
